# Generate Report for Handback
# Refresh the "generated" timestamps recorded on the handback-status report.
# These cells hold plain text timestamps (format "yyyy-MM-dd HH:mm:ss"),
# not real Excel date/time serials, so we assign them as strings.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the first file ---
# (column G, row 2). This timestamp is also mirrored as the "Correspond
# Handoff Datetime" for the same file on the de-de sheet (H2), since both
# workbook copies shared the same underlying string.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 23:04:19"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H2) and
#     "Correspond Handback DateTime" (K2) for the first file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-22 23:04:14"
$wsZhCn.Range("K2").Value = "2016-08-22 23:04:34"

# --- de-de sheet: "Correspond Handoff Datetime" (H2) and
#     "Correspond Handback DateTime" (K2) for the first file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-22 23:04:19"
$wsDeDe.Range("K2").Value = "2016-08-22 23:04:41"
